$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..Q (same set of values applies to every data row, 2..26)
$values = @(
    0.6408044419705359,    # B - r2
    -265.1326953808737,    # C - r2_sup
    -0.5545976795957159,   # D - r2_test
    0.7706479909155842,    # E - r2_val
    0.1701525824354797,    # F - r2_vt
    0.2132339996864685,    # G - mse
    157.9878643119225,     # H - mse_sup
    0.467468201747249,     # I - mse_test
    0.08282266875399483,   # J - mse_val
    0.2751454352506219,    # K - mse_vt
    0.2405142646481177,    # L - mape
    0.4617726710043249,    # M - rmse
    0.2163006006629874,    # N - r2_adj
    0.4814312896101858,    # O - rsd
    29.09073025240775,     # P - aic
    44.93611597569436      # Q - bic
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
